# Applies scheduled-runner profit/price updates to the Midgardsormr_Profits workbook
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, and WVR job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 778.125
$ws.Range("I18").Value = 778.125
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 778.125
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("N18").Value = -494.125
$ws.Range("H96").Value = 4464767.5
$ws.Range("I96").Value = 5952698.5
$ws.Range("J96").Value = 974.25
$ws.Range("K96").Value = 17858095.5
$ws.Range("L96").Value = 2922.75
$ws.Range("M96").Value = -17856722.5
$ws.Range("N96").Value = -5668.75
$ws.Range("H100").Value = 70159.664
$ws.Range("I100").Value = 70159.664
$ws.Range("K100").Value = 70159.664
$ws.Range("M100").Value = -69618.664
$ws.Range("H103").Value = 1951.9048
$ws.Range("I103").Value = 478.7
$ws.Range("J103").Value = 3291.182
$ws.Range("K103").Value = 1436.1
$ws.Range("L103").Value = 9873.545999999998
$ws.Range("M103").Value = -850.0999999999999
$ws.Range("N103").Value = -11045.546
$ws.Range("H112").Value = 6599.391
$ws.Range("J112").Value = 6712.711
$ws.Range("L112").Value = 20138.133
$ws.Range("N112").Value = -22354.133
$ws.Range("H116").Value = 2516701.2
$ws.Range("I116").Value = 5009700
$ws.Range("K116").Value = 5009700
$ws.Range("M116").Value = -5006258
$ws.Range("H137").Value = 19506.945
$ws.Range("J137").Value = 7666.6665
$ws.Range("L137").Value = 22999.9995
$ws.Range("N137").Value = -28099.9995
$ws.Range("H138").Value = 1607.9811
$ws.Range("J138").Value = 3297.889
$ws.Range("L138").Value = 9893.667000000001
$ws.Range("N138").Value = -20173.667
$ws.Range("H141").Value = 2289.4866
$ws.Range("I141").Value = 1797.8387
$ws.Range("K141").Value = 5393.5161
$ws.Range("M141").Value = -213.5160999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2162.7307
$ws.Range("J2").Value = 2705.4443
$ws.Range("L2").Value = 2705.4443
$ws.Range("N2").Value = -2931.4443
$ws.Range("H32").Value = 22574.744
$ws.Range("I32").Value = 17129.639
$ws.Range("J32").Value = 87916
$ws.Range("K32").Value = 17129.639
$ws.Range("L32").Value = 87916
$ws.Range("M32").Value = -16842.639
$ws.Range("N32").Value = -88490
$ws.Range("H61").Value = 8086.8887
$ws.Range("I61").Value = 1457.3334
$ws.Range("K61").Value = 1457.3334
$ws.Range("M61").Value = -1245.3334
$ws.Range("H102").Value = 9252.5
$ws.Range("I102").Value = 7336.6665
$ws.Range("K102").Value = 7336.6665
$ws.Range("M102").Value = -5714.6665
$ws.Range("H110").Value = 1657.25
$ws.Range("I110").Value = 1763.3334
$ws.Range("J110").Value = 1339
$ws.Range("K110").Value = 1763.3334
$ws.Range("L110").Value = 1339
$ws.Range("M110").Value = 281.6666
$ws.Range("N110").Value = -5429
$ws.Range("H116").Value = 2162.7307
$ws.Range("J116").Value = 2705.4443
$ws.Range("L116").Value = 2705.4443
$ws.Range("N116").Value = -7293.4443
$ws.Range("H136").Value = 8086.8887
$ws.Range("I136").Value = 1457.3334
$ws.Range("K136").Value = 4372.0002
$ws.Range("M136").Value = -1822.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2162.7307
$ws.Range("J3").Value = 2705.4443
$ws.Range("L3").Value = 2705.4443
$ws.Range("N3").Value = -2933.4443
$ws.Range("H105").Value = 4056.7273
$ws.Range("I105").Value = 4223.579
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 4223.579
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -2476.579
$ws.Range("N105").Value = -6494
$ws.Range("H107").Value = 17978.709
$ws.Range("I107").Value = 21111.76
$ws.Range("J107").Value = 4924.3335
$ws.Range("K107").Value = 21111.76
$ws.Range("L107").Value = 4924.3335
$ws.Range("M107").Value = -19191.76
$ws.Range("N107").Value = -8764.333500000001
$ws.Range("H134").Value = 2017.1744
$ws.Range("I134").Value = 1185.6025
$ws.Range("J134").Value = 10125
$ws.Range("K134").Value = 3556.8075
$ws.Range("L134").Value = 30375
$ws.Range("M134").Value = -1021.8075
$ws.Range("N134").Value = -35445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2780345
$ws.Range("I31").Value = 3127310
$ws.Range("K31").Value = 3127310
$ws.Range("M31").Value = -3127015
$ws.Range("H34").Value = 2780345
$ws.Range("I34").Value = 3127310
$ws.Range("K34").Value = 3127310
$ws.Range("M34").Value = -3127108
$ws.Range("H43").Value = 35552.332
$ws.Range("J43").Value = 35552.332
$ws.Range("L43").Value = 35552.332
$ws.Range("N43").Value = -35920.332
$ws.Range("H58").Value = 1159.091
$ws.Range("I58").Value = 1211.875
$ws.Range("J58").Value = 1018.3333
$ws.Range("K58").Value = 1211.875
$ws.Range("L58").Value = 1018.3333
$ws.Range("M58").Value = -1008.875
$ws.Range("N58").Value = -1424.3333
$ws.Range("H62").Value = 5630
$ws.Range("J62").Value = 7851.8184
$ws.Range("L62").Value = 7851.8184
$ws.Range("N62").Value = -9099.8184
$ws.Range("H65").Value = 5630
$ws.Range("J65").Value = 7851.8184
$ws.Range("L65").Value = 39259.092
$ws.Range("N65").Value = -45499.092
$ws.Range("H101").Value = 35552.332
$ws.Range("J101").Value = 35552.332
$ws.Range("L101").Value = 35552.332
$ws.Range("N101").Value = -42042.332
$ws.Range("H105").Value = 3309
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").Value = $null
$ws.Range("H136").Value = 1159.091
$ws.Range("I136").Value = 1211.875
$ws.Range("J136").Value = 1018.3333
$ws.Range("K136").Value = 3635.625
$ws.Range("L136").Value = 3054.9999
$ws.Range("M136").Value = -1085.625
$ws.Range("N136").Value = -8154.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 429.2
$ws.Range("I12").Value = 433.66666
$ws.Range("K12").Value = 1300.99998
$ws.Range("M12").Value = -1127.99998
$ws.Range("H87").Value = 10249.821
$ws.Range("I87").Value = 9606.429
$ws.Range("K87").Value = 28819.287
$ws.Range("M87").Value = -27571.287
$ws.Range("H90").Value = 10249.821
$ws.Range("I90").Value = 9606.429
$ws.Range("K90").Value = 86457.861
$ws.Range("M90").Value = -80217.861
$ws.Range("H140").Value = 3118.8333
$ws.Range("I140").Value = 3118.8333
$ws.Range("K140").Value = 9356.499899999999
$ws.Range("M140").Value = -4176.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 22004.4
$ws.Range("I102").Value = 25759.143
$ws.Range("K102").Value = 25759.143
$ws.Range("M102").Value = -24137.143
$ws.Range("H126").Value = 4815.923
$ws.Range("I126").Value = 2274.625
$ws.Range("J126").Value = 8882
$ws.Range("K126").Value = 6823.875
$ws.Range("L126").Value = 26646
$ws.Range("M126").Value = -4353.875
$ws.Range("N126").Value = -31586

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 26300.5
$ws.Range("J103").Value = 26300.5
$ws.Range("L103").Value = 26300.5
$ws.Range("N103").Value = -28644.5
$ws.Range("H104").Value = 8447.833000000001
$ws.Range("J104").Value = 8447.833000000001
$ws.Range("L104").Value = 8447.833000000001
$ws.Range("N104").Value = -15435.833
$ws.Range("H136").Value = 2515.3708
$ws.Range("I136").Value = 2244.8794
$ws.Range("K136").Value = 6734.638199999999
$ws.Range("M136").Value = -4184.638199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 8998
$ws.Range("I9").Value = 8998
$ws.Range("K9").Value = 8998
$ws.Range("M9").Value = -8858
$ws.Range("H57").Value = 83298
$ws.Range("J57").Value = 83298
$ws.Range("L57").Value = 83298
$ws.Range("N57").Value = -84806
$ws.Range("H82").Value = 1000
$ws.Range("I82").Value = 1000
$ws.Range("K82").Value = 1000
$ws.Range("M82").Value = -617
$ws.Range("H85").Value = 1000
$ws.Range("I85").Value = 1000
$ws.Range("K85").Value = 1000
$ws.Range("M85").Value = 326
$ws.Range("H100").Value = 1191.2142
$ws.Range("I100").Value = 368.2
$ws.Range("K100").Value = 736.4
$ws.Range("M100").Value = -195.4
$ws.Range("H122").Value = 16172895
$ws.Range("I122").Value = 17287244
$ws.Range("K122").Value = 51861732
$ws.Range("M122").Value = -51859282
$ws.Range("H136").Value = 8420.393
$ws.Range("I136").Value = 9541.392
$ws.Range("K136").Value = 28624.176
$ws.Range("M136").Value = -26074.176
